$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '39.458.65'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.67%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.161.31'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +3.17%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.81'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.45%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.623'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.06%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '64.44'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.88%  '

# Row 8
$ws.Range('E8').Value = '  +0.01%  '

# Row 9
$ws.Range('E9').Value = '  +2.75%  '

# Row 10
$ws.Range('E10').Value = '  +1.92%  '

# Row 11
$ws.Range('E11').Value = '  +0.55%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '16.07'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +4.73%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.482.43'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +3.21%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.28'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.75%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.814'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.53%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.56'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.41%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.165.27'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.72%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '39.429.16'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.76%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.83'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.06%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.13'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.92%  '

# Row 21
$ws.Range('E21').Value = '  +1.69%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '231.79'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.60%  '

# Row 23
$ws.Range('E23').Value = '  -0.01%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.51'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +5.83%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.36'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.92%  '

# Row 26
$ws.Range('E26').Value = '  +1.34%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '172.51'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.67%  '

# Row 28
$ws.Range('E28').Value = '  +1.97%  '

# Row 29
$ws.Range('E29').Value = '  +2.95%  '

# Row 30
$ws.Range('E30').Value = '  -0.50%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.69'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +8.77%  '

# Row 32
$ws.Range('E32').Value = '  +0.67%  '

# Row 33
$ws.Range('E33').Value = '  +2.35%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.76'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.15%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.10'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +8.87%  '

# Row 36
$ws.Range('E36').Value = '  -0.15%  '

# Row 37
$ws.Range('E37').Value = '  +0.40%  '

# Row 38
$ws.Range('E38').Value = '  -0.14%  '

# Row 39
$ws.Range('E39').Value = '  -0.02%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '104.13'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.79%  '

# Row 41
$ws.Range('E41').Value = '  +0.94%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '17.92'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.19%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.539.68'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.28%  '

# Row 44
$ws.Range('E44').Value = '  +3.90%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.93'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +3.08%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0926'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.83%  '

# Row 47
$ws.Range('E47').Value = '  +0.50%  '

# Row 48
$ws.Range('E48').Value = '  +5.86%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.20'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.39%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.364.75'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.11%  '

# Row 51
$ws.Range('E51').Value = '  +0.00%  '
